# Weekly update: insert a new record for the latest week at row 16,
# pushing the existing rows 16-39 down to 17-40.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(16).Insert()

$ws.Range("A16").Value = 11
$ws.Range("B16").Value = "Vega Monumental Concepción"
$ws.Range("C16").Value = "Bíobío"
$ws.Range("D16").Value = 44810
$ws.Range("E16").Value = 8
$ws.Range("F16").Value = 100114007
$ws.Range("G16").Value = "Jengibre"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 50
$ws.Range("K16").Value = 11000
$ws.Range("L16").Value = 12000
$ws.Range("M16").Value = 11600
$ws.Range("N16").Value = "$/caja 13 kilos"
$ws.Range("O16").Value = "Perú"
$ws.Range("P16").Value = 892
$ws.Range("Q16").Value = 13
$ws.Range("R16").Value = "Hortaliza"
